$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C8 value from 0 to 111
$ws.Range("C8").Value = 111

# Collapse all defined columns (A through the last used, plus the default rest)
$ws.Columns.Item(1).EntireColumn.OutlineLevel = $ws.Columns.Item(1).EntireColumn.OutlineLevel
for ($i = 1; $i -le 12; $i++) {
    $ws.Columns.Item($i).EntireColumn.Collapsed = $true
}
